# "Resolução exercício aula 3"
# Updates the test-scenario tracking sheet: several Status (column H)
# values are corrected, row 16's Criticidade (column G) is filled in,
# row 19's Resultado Esperado (column F) is corrected, and the
# workbook selection is moved to H19.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Status (column H) corrections ---------------------------------
$ws.Range("H2").Value  = "Não iniciada"
$ws.Range("H6").Value  = "Não iniciada"
$ws.Range("H9").Value  = "Não iniciada"
$ws.Range("H11").Value = "Sucesso"
$ws.Range("H12").Value = "Sucesso"
$ws.Range("H14").Value = "Sucesso"
$ws.Range("H15").Value = "Sucesso"
$ws.Range("H16").Value = "Sucesso"
$ws.Range("H18").Value = "Falha"
$ws.Range("H19").Value = "Sucesso"
$ws.Range("H24").Value = "Não iniciada"

# --- Row 16: fill in the missing Criticidade -------------------------
$ws.Range("G16").Value = "Média"

# --- Row 19: correct the expected result text -----------------------
$ws.Range("F19").Value = "Então o sistema exibe uma mensagem de erro"

# --- Move the active selection to H19 --------------------------------
$ws.Range("H19").Select()
